# Update the division-problem table: replace each cell's expression with
# its new value. Cell text is addressed positionally (row/column) rather
# than via global Find/Replace because several source/target strings repeat
# or collide (e.g. "24÷9=" appears twice with different replacements, and
# "66÷4=" is both a source and a target elsewhere), so a blind text search
# could match the wrong occurrence.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of (row, column) -> new expression, in document order.
$updates = @(
    @{Row=1;  Col=1; Text="75÷4="},
    @{Row=1;  Col=2; Text="45÷5="},
    @{Row=1;  Col=3; Text="10÷5="},
    @{Row=1;  Col=4; Text="55÷6="},
    @{Row=1;  Col=5; Text="14÷8="},

    @{Row=5;  Col=1; Text="40÷5="},
    @{Row=5;  Col=2; Text="52÷7="},
    @{Row=5;  Col=3; Text="30÷5="},
    @{Row=5;  Col=4; Text="16÷4="},
    @{Row=5;  Col=5; Text="59÷9="},

    @{Row=9;  Col=1; Text="74÷8="},
    @{Row=9;  Col=2; Text="64÷4="},
    @{Row=9;  Col=3; Text="38÷9="},
    @{Row=9;  Col=4; Text="39÷3="},
    @{Row=9;  Col=5; Text="96÷4="},

    @{Row=13; Col=1; Text="66÷4="},
    @{Row=13; Col=2; Text="90÷7="},
    @{Row=13; Col=3; Text="10÷7="},
    @{Row=13; Col=4; Text="65÷3="},
    @{Row=13; Col=5; Text="18÷2="},

    @{Row=17; Col=1; Text="82÷9="},
    @{Row=17; Col=2; Text="20÷4="},
    @{Row=17; Col=3; Text="59÷6="},
    @{Row=17; Col=4; Text="48÷9="},
    @{Row=17; Col=5; Text="25÷5="}
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
